# 3.1. Hiển thị thông báo sau khi đổi trạng thái
# Adds a new row (row 8) documenting "2. Tính năng thay đổi thứ tự sản phẩm"
# with its commit link, and puts a thin box border around the whole
# A2:C8 data block (excluding the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- add the new data row -------------------------------------------------
$ws.Range("B8").Value = "2. Tính năng thay đổi thứ tự sản phẩm"
$ws.Range("C8").Value = "https://github.com/nguyentienminh07102004/product-management/commit/a36764f2bc798b74ffb3bdda7168d990212f567d"
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/nguyentienminh07102004/product-management/commit/a36764f2bc798b74ffb3bdda7168d990212f567d") | Out-Null

# --- apply a thin box border around every data cell (rows 2-8) -----------
$rng = $ws.Range("A2:C8")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# --- tidy up the view: clear the old scroll/selection, select B12 --------
$ws.Activate()
$ws.Range("B12").Select() | Out-Null

Write-Host "done"
